$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D42").Value = "오류 처리 관련 포스팅 정리"
$ws.Range("E42").Value = "https://kjk92.tistory.com/95"

$ws.Range("D52").Value = "Relative Risk Regression"

$wb.Save()
